$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# --- Data update (bot daily refresh) ---------------------------------
# 2020-04-25 (row 91) numbers have now arrived.
$ws.Cells.Item(91, 2).Value = 524
$ws.Cells.Item(91, 3).Value = 29514
$ws.Cells.Item(91, 4).Value = 0
$ws.Cells.Item(91, 5).Value = 6458

# Insert a new row ahead of the footnote row (currently row 92) so the
# note shifts down to row 93. The inserted row inherits the formatting
# of the row above it (row 91), matching the other data rows' styling.
$ws.Rows.Item(92).Insert()
$ws.Cells.Item(92, 1).Value = 43947

# --- View state ---------------------------------------------------------
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
$ws.Range("B92").Select()

# --- Print area -----------------------------------------------------
# (Set via the defined-name object directly rather than
# $ws.PageSetup.PrintArea so the sheet qualifier is written unquoted,
# matching the workbook's existing definedName formatting.)
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$95"
    }
}
